# Refresh the cryptos table (coinranking.com snapshot) for the
# GitHub Actions scheduled update: prices + 1h volume % per coin,
# row 5/6 also swap because BNB and Solana's rank order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "49.871.15"
$ws.Range("E2").Value = "  +3.67%  "

# Row 3
$ws.Range("D3").Value = "2.647.19"
$ws.Range("E3").Value = "  +5.80%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'327.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'111.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.50%  "

# Row 7
$ws.Range("E7").Value = "  +1.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.559"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.23%  "

# Row 10
$ws.Range("D10").Value = "'40.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "

# Row 11
$ws.Range("D11").Value = "'20.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.74%  "

# Row 12
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("E13").Value = "  +0.76%  "

# Row 14
$ws.Range("E14").Value = "  +3.20%  "

# Row 15
$ws.Range("D15").Value = "3.062.27"
$ws.Range("E15").Value = "  +5.85%  "

# Row 16
$ws.Range("D16").Value = "2.666.33"
$ws.Range("E16").Value = "  +6.49%  "

# Row 17
$ws.Range("E17").Value = "  +5.47%  "

# Row 18
$ws.Range("D18").Value = "49.824.81"
$ws.Range("E18").Value = "  +3.85%  "

# Row 19
$ws.Range("D19").Value = "'13.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.26%  "

# Row 20
$ws.Range("E20").Value = "  +2.06%  "

# Row 21
$ws.Range("D21").Value = "'2.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.64%  "

# Row 22
$ws.Range("E22").Value = "  +2.13%  "

# Row 23
$ws.Range("D23").Value = "'72.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "

# Row 24
$ws.Range("D24").Value = "'280.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "

# Row 25
$ws.Range("D25").Value = "'2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.29%  "

# Row 26
$ws.Range("D26").Value = "'27.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.43%  "

# Row 27
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("E28").Value = "  -2.36%  "

# Row 29
$ws.Range("E29").Value = "  +4.71%  "

# Row 30
$ws.Range("D30").Value = "'9.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.56%  "

# Row 31
$ws.Range("E31").Value = "  +1.71%  "

# Row 32
$ws.Range("D32").Value = "'49.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").Value = "'19.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.71%  "

# Row 34
$ws.Range("E34").Value = "  +2.58%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").Value = "'0.0798"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.88%  "

# Row 37
$ws.Range("E37").Value = "  +6.34%  "

# Row 38
$ws.Range("D38").Value = "'4.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "

# Row 39
$ws.Range("E39").Value = "  +8.06%  "

# Row 40
$ws.Range("D40").Value = "'126.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.51%  "

# Row 41
$ws.Range("E41").Value = "  +1.66%  "

# Row 42
$ws.Range("E42").Value = "  +1.76%  "

# Row 43
$ws.Range("D43").Value = "'22.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.66%  "

# Row 44
$ws.Range("E44").Value = "  +3.71%  "

# Row 45
$ws.Range("E45").Value = "  +7.88%  "

# Row 46
$ws.Range("D46").Value = "2.067.18"
$ws.Range("E46").Value = "  +2.31%  "

# Row 47
$ws.Range("D47").Value = "'2.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.62%  "

# Row 48
$ws.Range("D48").Value = "'1.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.03%  "

# Row 49
$ws.Range("E49").Value = "  +0.93%  "

# Row 50
$ws.Range("D50").Value = "'5.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.54%  "

# Row 51
$ws.Range("D51").Value = "'82.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.04%  "
